{"js": "// Fix the typo \"raining\" -> \"training\" in the \"Besides described above, ...\"\n// paragraph (insert a single \"t\" right before \"raining at\"), and relocate the\n// \"_GoBack\" bookmark from the blank paragraph above onto this paragraph,\n// right after the sentence and before its trailing line break.\n\nconst body = context.document.body;\n\n// 1) Insert the missing \"t\" so \"Conducting raining at\" becomes\n//    \"Conducting training at\".\nconst typoResults = body.search(\"raining at\", { matchCase: false });\ntypoResults.load(\"items\");\nawait context.sync();\n\nif (typoResults.items.length > 0) {\n  const typoRange = typoResults.items[0];\n  const insertionPoint = typoRange.getRange(Word.RangeLocation.start);\n  insertionPoint.insertText(\"t\", Word.InsertLocation.before);\n  await context.sync();\n}\n\n// 2) Move the \"_GoBack\" bookmark: delete it from wherever it currently is\n//    (the empty paragraph right above), then re-insert it right after\n//    \"...Data Engineering & Analytics\" (and before the line break that\n//    follows it) in the training paragraph.\ncontext.document.deleteBookmark(\"_GoBack\");\n\nconst anchorResults = body.search(\"Data Engineering & Analytics\", { matchCase: false });\nanchorResults.load(\"items\");\nawait context.sync();\n\nif (anchorResults.items.length > 0) {\n  const anchorRange = anchorResults.items[0];\n  const bookmarkPoint = anchorRange.getRange(Word.RangeLocation.end);\n  bookmarkPoint.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Fix the typo \"raining\" -> \"training\" in the \"Besides described above, ...\"\n# paragraph (insert a single \"t\" right before \"raining at\"), and relocate the\n# \"_GoBack\" bookmark from the blank paragraph above onto this paragraph,\n# right after the sentence and before its trailing line break.\n\n$d = $word.ActiveDocument\n\n# 1) Insert the missing \"t\" so \"Conducting raining at\" becomes\n#    \"Conducting training at\".\n$findRng = $d.Content\n$findRng.Find.ClearFormatting()\n$findRng.Find.Text = \"raining at\"\n$found = $findRng.Find.Execute()\nif ($found) {\n  $insertionPoint = $d.Range($findRng.Start, $findRng.Start)\n  $insertionPoint.InsertBefore(\"t\")\n}\n\n# 2) Move the \"_GoBack\" bookmark: delete it from wherever it currently is\n#    (the empty paragraph right above), then re-insert it right after\n#    \"...Data Engineering & Analytics\" (and before the line break that\n#    follows it) in the training paragraph.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n  $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$anchorRng = $d.Content\n$anchorRng.Find.ClearFormatting()\n$anchorRng.Find.Text = \"Data Engineering & Analytics\"\n$anchorFound = $anchorRng.Find.Execute()\nif ($anchorFound) {\n  $bookmarkPoint = $d.Range($anchorRng.End, $anchorRng.End)\n  $d.Bookmarks.Add(\"_GoBack\", $bookmarkPoint)\n}\n"}
